$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F666").Value = "[73.5]"
$ws.Range("G666").Value = 73.5

$ws.Range("F667").Value = "[81.5]"
$ws.Range("G667").Value = 81.5
$ws.Range("H667").Value = 81.5

$ws.Range("F668").Value = "[67.375]"
$ws.Range("G668").Value = 67.375
$ws.Range("H668").Value = 81.5

$ws.Range("F669").Value = "[19.0]"
$ws.Range("G669").Value = 19
$ws.Range("H669").Value = 81.5

$ws.Range("F670").Value = "[21.0]"
$ws.Range("G670").Value = 21
$ws.Range("H670").Value = 81.5

$ws.Range("H671").Value = 81.5

$ws.Range("F672").Value = "[85.0]"
$ws.Range("G672").Value = 85
$ws.Range("H672").Value = 85

$ws.Range("F673").Value = "[69.5]"
$ws.Range("G673").Value = 69.5

$ws.Range("F674").Value = "[76.625]"
$ws.Range("G674").Value = 76.625

$ws.Range("F675").Value = "[81.5]"
$ws.Range("G675").Value = 81.5

$ws.Range("F690").Value = "[75.26333333333334]"
$ws.Range("G690").Value = 75.26333333333334

$ws.Range("F691").Value = "[75.15333333333332]"
$ws.Range("G691").Value = 75.15333333333332

$ws.Range("F710").Value = "[73.5]"
$ws.Range("G710").Value = 73.5

$ws.Range("F711").Value = "[81.5]"
$ws.Range("G711").Value = 81.5
$ws.Range("H711").Value = 81.5

$ws.Range("F712").Value = "[67.375]"
$ws.Range("G712").Value = 67.375
$ws.Range("H712").Value = 81.5

$ws.Range("F713").Value = "[19.0]"
$ws.Range("G713").Value = 19
$ws.Range("H713").Value = 81.5

$ws.Range("F714").Value = "[21.0]"
$ws.Range("G714").Value = 21
$ws.Range("H714").Value = 81.5

$ws.Range("F715").Value = "[75.26333333333334]"
$ws.Range("G715").Value = 75.26333333333334
$ws.Range("H715").Value = 81.5

$ws.Range("F716").Value = "[76.625]"
$ws.Range("G716").Value = 76.625
$ws.Range("H716").Value = 81.5

$ws.Range("F717").Value = "[69.5]"
$ws.Range("G717").Value = 69.5
$ws.Range("H717").Value = 81.5

$ws.Range("F718").Value = "[85.0]"
$ws.Range("G718").Value = 85

$ws.Range("F751").Value = "[19.0]"
$ws.Range("G751").Value = 19

$ws.Range("F754").Value = "[59.125]"
$ws.Range("G754").Value = 59.125

$ws.Range("F755").Value = "[75.26333333333334]"
$ws.Range("G755").Value = 75.26333333333334

$ws.Range("F756").Value = "[75.15333333333332]"
$ws.Range("G756").Value = 75.15333333333332

$ws.Range("F757").Value = "[21.0]"
$ws.Range("G757").Value = 21

$ws.Range("F758").Value = "[56.625]"
$ws.Range("G758").Value = 56.625

$ws.Range("F759").Value = "[74.08666666666666]"
$ws.Range("G759").Value = 74.08666666666666

$ws.Range("F760").Value = "[85.0]"
$ws.Range("G760").Value = 85
$ws.Range("H760").Value = 85

$ws.Range("F761").Value = "[69.5]"
$ws.Range("G761").Value = 69.5
$ws.Range("H761").Value = 85

$ws.Range("F762").Value = "[76.625]"
$ws.Range("G762").Value = 76.625

$ws.Range("F800").Value = "[75.26333333333334]"
$ws.Range("G800").Value = 75.26333333333334

$ws.Range("F801").Value = "[75.15333333333332]"
$ws.Range("G801").Value = 75.15333333333332

$ws.Range("F845").Value = "[81.5]"
$ws.Range("G845").Value = 81.5
$ws.Range("H845").Value = 83

$ws.Range("F846").Value = "[81.5]"
$ws.Range("G846").Value = 81.5
$ws.Range("H846").Value = 83

$ws.Range("F847").Value = "[74.125]"
$ws.Range("G847").Value = 74.125
$ws.Range("H847").Value = 83

$ws.Range("F848").Value = "[85.0]"
$ws.Range("G848").Value = 85

$ws.Range("F849").Value = "[21.0]"
$ws.Range("G849").Value = 21

$ws.Range("F850").Value = "[19.0]"
$ws.Range("G850").Value = 19

$ws.Range("F851").Value = "[73.5]"
$ws.Range("G851").Value = 73.5

$ws.Range("F852").Value = "[79.0]"
$ws.Range("G852").Value = 79

$ws.Range("F853").Value = "[75.15333333333332]"
$ws.Range("G853").Value = 75.15333333333332

$ws.Range("F854").Value = "[76.625]"
$ws.Range("G854").Value = 76.625

$ws.Range("F855").Value = "[67.375]"
$ws.Range("G855").Value = 67.375

$ws.Range("F856").Value = "[77.0]"
$ws.Range("G856").Value = 77
